# Auto-generated edit script applying price/coin updates per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '242.51'

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '21.51'

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '5.224'

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '0.05601'

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '3.369'

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '6.371'

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.9368'

$ws.Cells.Item(10, 2).Value = 'WazirX'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.1424'
$ws.Cells.Item(10, 5).Value = '9WazirXWRX'

$ws.Cells.Item(11, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.07301'
$ws.Cells.Item(11, 5).Value = '10MandalaExchangeTokenMDX'

$ws.Cells.Item(12, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.03129'
$ws.Cells.Item(12, 5).Value = '11LiechtensteinCryptoassetsExchangeLCX'

$ws.Cells.Item(13, 2).Value = 'BitrueCoin'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '0.03024'
$ws.Cells.Item(13, 5).Value = '12BitrueCoinBTR'

$ws.Cells.Item(14, 2).Value = 'BitMartToken'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '0.09274'
$ws.Cells.Item(14, 5).Value = '13BitMartTokenBMX'

$ws.Cells.Item(15, 2).Value = 'MCDex'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '3.611'
$ws.Cells.Item(15, 5).Value = '14MCDexMCB'

$ws.Cells.Item(16, 2).Value = 'BitForexToken'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '0.001649'
$ws.Cells.Item(16, 5).Value = '15BitForexTokenBF'

$ws.Cells.Item(17, 2).Value = 'CoinExToken'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '0.04693'
$ws.Cells.Item(17, 5).Value = '16CoinExTokenCET'

$ws.Cells.Item(18, 2).Value = 'One'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '0.0005800'
$ws.Cells.Item(18, 5).Value = '17OneONE'

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '0.006356'

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '0.004981'

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '0.0001500'

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '0.0003099'

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '3.763'

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '2.094'

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '0.3268'

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '0.03918'

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '0.006875'

$ws.Cells.Item(42, 2).Value = 'BKEXToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '0.1035'
$ws.Cells.Item(42, 5).Value = '41BKEXTokenBKK'

$ws.Cells.Item(43, 2).Value = 'CEJI'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '0.002921'
$ws.Cells.Item(43, 5).Value = '42CEJICEJI'

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.008301'

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.00005945'

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '0.0005500'

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '0.6824'

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '0.06732'

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.00002100'
